$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains numeric-looking text (prices) that must stay as TEXT,
# matching the source workbook (t="inlineStr"/shared-string cells, no numeric
# coercion). Temporarily force the column to Text format while writing the
# new values, then restore the original style/format so no stray styles leak
# into the saved workbook.
$dRange = $ws.Range("D2:D51")
$origStyle = $dRange.Style()
$dRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "58.632.64"
$ws.Range("E2").Value = "  -5.04%  "

# Row 3
$ws.Range("D3").Value = "2.468.36"
$ws.Range("E3").Value = "  -4.36%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").Value = "534.80"
$ws.Range("E5").Value = "  -3.41%  "

# Row 6
$ws.Range("D6").Value = "144.28"
$ws.Range("E6").Value = "  -6.29%  "

# Row 7
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.39%  "

# Row 8
$ws.Range("E8").Value = "  -4.82%  "

# Row 9
$ws.Range("D9").Value = "2.492.49"
$ws.Range("E9").Value = "  -3.61%  "

# Row 10
$ws.Range("E10").Value = "  -4.38%  "

# Row 11
$ws.Range("E11").Value = "  -2.26%  "

# Row 12
$ws.Range("D12").Value = "5.61"
$ws.Range("E12").Value = "  +2.42%  "

# Row 13
$ws.Range("E13").Value = "  -3.07%  "

# Row 14
$ws.Range("D14").Value = "2.898.33"
$ws.Range("E14").Value = "  -4.62%  "

# Row 15
$ws.Range("D15").Value = "23.76"
$ws.Range("E15").Value = "  -6.50%  "

# Row 16
$ws.Range("D16").Value = "58.536.15"
$ws.Range("E16").Value = "  -5.06%  "

# Row 17
$ws.Range("E17").Value = "  -4.10%  "

# Row 18
$ws.Range("D18").Value = "2.478.38"
$ws.Range("E18").Value = "  -4.15%  "

# Row 19
$ws.Range("D19").Value = "11.31"
$ws.Range("E19").Value = "  -2.21%  "

# Row 20
$ws.Range("D20").Value = "4.32"
$ws.Range("E20").Value = "  -4.58%  "

# Row 21
$ws.Range("D21").Value = "323.50"
$ws.Range("E21").Value = "  -4.41%  "

# Row 22
$ws.Range("D22").Value = "0.996"
$ws.Range("E22").Value = "  -0.18%  "

# Row 23
$ws.Range("D23").Value = "5.74"
$ws.Range("E23").Value = "  -4.94%  "

# Row 24
$ws.Range("D24").Value = "60.48"
$ws.Range("E24").Value = "  -3.87%  "

# Row 25
$ws.Range("E25").Value = "  -11.26%  "

# Row 26
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  -0.11%  "

# Row 27
$ws.Range("E27").Value = "  -4.46%  "

# Row 28
$ws.Range("D28").Value = "2.579.55"
$ws.Range("E28").Value = "  -4.60%  "

# Row 29
$ws.Range("D29").Value = "7.73"
$ws.Range("E29").Value = "  -4.17%  "

# Row 30
$ws.Range("E30").Value = "  -0.77%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0774"
$ws.Range("E31").Value = "  -7.39%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.79"
$ws.Range("E32").Value = "  -6.62%  "

# Row 33
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "1.23"
$ws.Range("E33").Value = "  -5.37%  "

# Row 34
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  -0.30%  "

# Row 35
$ws.Range("D35").Value = "157.97"
$ws.Range("E35").Value = "  -1.46%  "

# Row 36
$ws.Range("E36").Value = "  -0.84%  "

# Row 37
$ws.Range("E37").Value = "  -3.79%  "

# Row 38
$ws.Range("E38").Value = "  -5.87%  "

# Row 39
$ws.Range("D39").Value = "1.61"
$ws.Range("E39").Value = "  -10.13%  "

# Row 40
$ws.Range("E40").Value = "  -4.80%  "

# Row 41
$ws.Range("D41").Value = "306.54"
$ws.Range("E41").Value = "  -8.88%  "

# Row 42
$ws.Range("D42").Value = "36.50"
$ws.Range("E42").Value = "  -2.49%  "

# Row 43
$ws.Range("D43").Value = "3.71"
$ws.Range("E43").Value = "  -5.39%  "

# Row 44
$ws.Range("D44").Value = "0.813"
$ws.Range("E44").Value = "  -8.86%  "

# Row 45
$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  -0.09%  "

# Row 46
$ws.Range("D46").Value = "10.76"
$ws.Range("E46").Value = "  -1.61%  "

# Row 47
$ws.Range("D47").Value = "0.592"
$ws.Range("E47").Value = "  -2.32%  "

# Row 48
$ws.Range("D48").Value = "124.35"
$ws.Range("E48").Value = "  -0.02%  "

# Row 49
$ws.Range("E49").Value = "  -4.13%  "

# Row 50
$ws.Range("D50").Value = "0.0519"
$ws.Range("E50").Value = "  -4.70%  "

# Row 51
$ws.Range("D51").Value = "0.0228"
$ws.Range("E51").Value = "  -4.89%  "

# Restore original style/number format for column D
$dRange.Style = $origStyle

Write-Host "Update complete"